$d = $word.ActiveDocument

function Get-ParagraphStartingWith($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text.StartsWith($needle)) {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Paragraph: "El sistema deberá de ser usado solo por el personal..."
# Prepend a highlighted label run "No funcional de seguridad" followed
# by a space, keeping the original sentence (now prefixed with a
# space) in its own, unhighlighted run.
# ---------------------------------------------------------------------
$p1 = Get-ParagraphStartingWith $d "El sistema deber"
$insertPoint1 = $p1.Range.Duplicate
$insertPoint1.Collapse(1)   # wdCollapseStart

$label1 = "No funcional de seguridad"
$insertPoint1.InsertBefore($label1 + " ")

$labelStart1 = $p1.Range.Start
$labelRange1 = $d.Range($labelStart1, $labelStart1 + $label1.Length)
$labelRange1.Font.HighlightColorIndex = 6   # wdRed

# ---------------------------------------------------------------------
# Paragraph: "Se pide que el sistema porte los colores de la
# institución, así como su logotipo."
# Make the whole paragraph italic, then prepend a highlighted label run
# "Requisito de interfaz" followed by a space, keeping the original
# sentence (now prefixed with a space) in its own run (italic, no
# highlight).
# ---------------------------------------------------------------------
$p5 = Get-ParagraphStartingWith $d "Se pide que el sistema porte"
$p5.Range.Font.Italic = 1

$insertPoint5 = $p5.Range.Duplicate
$insertPoint5.Collapse(1)   # wdCollapseStart

$label5 = "Requisito de interfaz"
$insertPoint5.InsertBefore($label5 + " ")

$labelStart5 = $p5.Range.Start
$labelRange5 = $d.Range($labelStart5, $labelStart5 + $label5.Length)
$labelRange5.Font.Italic = 1
$labelRange5.Font.HighlightColorIndex = 6   # wdRed

Write-Output "Applied business-rule labels."
